$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("K3").Value = '("success factors" OR "factors" OR "aspects" OR "items" OR "elements" OR "drivers" OR "motivators" OR "variables") AND ("Outsourcing" OR "global software development" OR "geographically distributed development" OR "offshore development" OR "multisite development" OR "collaborative software development") AND ("IaaS" OR "PaaS" OR "SaaS" OR "XaaS" OR "Infrastructure as a Service" OR "Platform as a Service" OR "Software as a Service" OR "IT service" OR "Application Service" OR "ASP") AND ("cloud computing" OR "cloud platform" OR "cloud provider" OR "cloud service" OR "cloud offering")'
$ws.Range("Q3").Value = "['application service', 'asp', 'aspects', 'cloud computing', 'cloud offering', 'cloud platform', 'cloud provider', 'cloud service', 'collaborative software development', 'drivers', 'elements', 'factors', 'geographically distributed development', 'global software development', 'iaas', 'infrastructure as a service', 'it service', 'items', 'motivators', 'multisite development', 'offshore development', 'outsourcing', 'paas', 'platform as a service', 'saas', 'software as a service', 'success factors', 'variables', 'xaas']"
$ws.Range("S3").Value = 0.02631578947368421
$ws.Range("T3").Value = 0.04166666666666666

# Row 5
$ws.Range("K5").Value = '(blockchain OR ((distributed OR decentralized) AND (ledger OR platform OR "autonomous organization"))) AND (governance OR management OR ecosystem)'
$ws.Range("Q5").Value = "['autonomous organization', 'blockchain', 'decentralized', 'distributed', 'ecosystem', 'governance', 'ledger', 'management', 'platform']"
$ws.Range("S5").Value = 0.1
$ws.Range("T5").Value = 0.1290322580645161

# Row 15
$ws.Range("K15").Value = '(smell OR “design flaw” OR disharmony OR “code anomaly” OR “design anomaly” OR anti-pattern) AND (experiment OR empirical OR survey OR ethnography OR “action research” OR “exploratory analysis” OR study OR controlled)'
$ws.Range("Q15").Value = "['anti-pattern', 'controlled', 'disharmony', 'empirical', 'ethnography', 'experiment', 'smell', 'study', 'survey', '“action research”', '“code anomaly”', '“design anomaly”', '“design flaw”', '“exploratory analysis”']"
$ws.Range("T15").Value = 0.04166666666666666

# Row 17
$ws.Range("K17").Value = '(agile AND practice AND select) OR (agile AND method AND tailoring) OR (scrum AND practice AND adoption) OR (scrum AND tailoring) OR (scrum AND practice AND select) OR (xp AND practice AND adoption) OR (xp AND tailoring) OR (xp AND practice AND select) OR (kanban AND practice AND adoption) OR (kanban AND tailoring) OR (kanban AND practice AND select) OR (lean AND practice AND adoption) OR (lean AND tailoring) OR (lean AND practice AND select) OR (fdd AND practice AND adoption) OR (fdd AND tailoring) OR (fdd AND practice AND select) OR (“feature driven development” AND practice AND adoption) OR (“feature driven development” AND tailoring) OR (“feature driven development” AND practice AND select)'
$ws.Range("Q17").Value = "['adoption', 'adoption', 'adoption', 'adoption', 'adoption', 'adoption', 'agile', 'agile', 'fdd', 'fdd', 'fdd', 'kanban', 'kanban', 'kanban', 'lean', 'lean', 'lean', 'method', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'practice', 'scrum', 'scrum', 'scrum', 'select', 'select', 'select', 'select', 'select', 'select', 'select', 'tailoring', 'tailoring', 'tailoring', 'tailoring', 'tailoring', 'tailoring', 'tailoring', 'xp', 'xp', 'xp', '“feature driven development”', '“feature driven development”', '“feature driven development”']"
$ws.Range("S17").Value = 0.06666666666666667
$ws.Range("T17").Value = 0.1282051282051282

# Row 25
$ws.Range("K25").Value = '(“intention to leave” or “intention for turnover” or “intention to quit” or “intention to withdraw” or “intention to stay” or “turnover intention” or “leave intention” or “quit intention” or “stay intention” or “withdrawal intention” or “employee retention” or “employees retention” or “personnel retention” or “worker retention” or “workers retention” or “manager retention” or “managers retention” or “professional retention” or “professionals retention”) AND (“IS worker” or “IT worker” or “IS employee” or “IT employee” or “IS personnel” or “IT personnel” or “IS manager” or “IT manager” or “IS professional” or “IT professional” or “IS workforce” or “IT workforce” or “Software engineer” or “software developer” or “software programmer” or “software manager” or “system analyst” or “software designer” or “software project manager”)'
$ws.Range("Q25").Value = "['“employee retention”', '“employees retention”', '“intention for turnover”', '“intention to leave”', '“intention to quit”', '“intention to stay”', '“intention to withdraw”', '“is employee”', '“is manager”', '“is personnel”', '“is professional”', '“is worker”', '“is workforce”', '“it employee”', '“it manager”', '“it personnel”', '“it professional”', '“it worker”', '“it workforce”', '“leave intention”', '“manager retention”', '“managers retention”', '“personnel retention”', '“professional retention”', '“professionals retention”', '“quit intention”', '“software designer”', '“software developer”', '“software engineer”', '“software manager”', '“software programmer”', '“software project manager”', '“stay intention”', '“system analyst”', '“turnover intention”', '“withdrawal intention”', '“worker retention”', '“workers retention”']"

# Row 26
$ws.Range("K26").Value = '(agile OR kanban OR scrum OR lean OR “extreme programming” OR “design thinking”) AND (“user experience” OR ux OR usability OR hcd OR hci OR hmi OR ucd)'
$ws.Range("Q26").Value = "['agile', 'hcd', 'hci', 'hmi', 'kanban', 'lean', 'scrum', 'ucd', 'usability', 'ux', '“design thinking”', '“extreme programming”', '“user experience”']"
$ws.Range("S26").Value = 0.1379310344827586
$ws.Range("T26").Value = 0.1081081081081081
